# The sheet contains a daily price log for "Perejil" (Parsley) at the
# Mercado Mayorista Lo Valledor de Santiago market, with one data row
# per day starting at row 2. A new day's record needs to be inserted
# in the middle of the existing data (at row 473), which pushes all
# subsequent rows (old 473:578) down by one (to 474:579).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 473; this shifts the previous
# rows 473:578 down to 474:579 and extends the used range to row 579.
$ws.Rows("473:473").Insert()

# Populate the newly inserted row 473 with the new daily record.
$ws.Range("A473").Value = 6
$ws.Range("B473").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C473").Value = "Metropolitana"
$ws.Range("D473").Value = 44798
$ws.Range("E473").Value = 13
$ws.Range("F473").Value = 100112044
$ws.Range("G473").Value = "Perejil"
$ws.Range("H473").Value = "Sin especificar"
$ws.Range("I473").Value = "Primera"
$ws.Range("J473").Value = 190
$ws.Range("K473").Value = 18000
$ws.Range("L473").Value = 19000
$ws.Range("M473").Value = 18368
$ws.Range("N473").Value = "`$/docena de atados"
$ws.Range("O473").Value = "Región Metropolitana"
$ws.Range("P473").Value = 6123
$ws.Range("Q473").Value = 3
$ws.Range("R473").Value = "Hortaliza"
